# Migrate to the newest version of Closed.XML
#
# The upgraded ClosedXML writer no longer emits degenerate single-cell
# "merges" (e.g. <mergeCell ref="B3:B3"/>) and fixes a couple of cells that
# were previously serialized as blank/placeholder values inside merged
# ranges. Reproduce the equivalent final worksheet state here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the bogus single-cell "merges" left over from the old writer.
#    These ranges were never really merged (MergeArea == the cell itself)
#    so un-merging them just removes the stray <mergeCell> entries.
# ---------------------------------------------------------------------
$bogusSingleCellMerges = @("B3","B6","D2","D4","D5","E5","E6","F2","F5","G2","G4","G5","H2")
foreach ($addr in $bogusSingleCellMerges) {
    $ws.Range($addr).UnMerge()
}

# ---------------------------------------------------------------------
# 2. Fix D6: it is the secondary cell of the C6:D6 merge and should end up
#    completely empty (no type, no value) instead of the old placeholder
#    blank-string cell. Un-merge so the cell is directly addressable,
#    clear it, then restore the merge.
# ---------------------------------------------------------------------
$ws.Range("C6:D6").UnMerge()
$ws.Range("D6").ClearContents()
$ws.Range("C6:D6").Merge()

# ---------------------------------------------------------------------
# 3. Fix I6: it was the secondary cell of the H6:I6 merge holding an empty
#    placeholder string. The new output un-merges H6:I6 entirely and
#    stores real text content ("21.02.2018 0:00:00") in I6, formatted with
#    the default (General) style instead of the inherited date format.
# ---------------------------------------------------------------------
$ws.Range("H6:I6").UnMerge()
$ws.Range("I6").Value = "21.02.2018 0:00:00"
$ws.Range("B2").Copy()
$ws.Range("I6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. E2 and F2 were previously two separate degenerate single-cell
#    "merges"; the new writer represents this as one real E2:F2 merge.
# ---------------------------------------------------------------------
$ws.Range("E2:F2").Merge()
